# feat: add 2022-Q4 data
#
# 1. Insert a brand-new worksheet "2022-Q4" right after "总计" (before the
#    existing "2022-Q3" sheet), built from a copy of the "2022-Q3" sheet so
#    it inherits the exact same layout/styles, then overwritten with the
#    2022-Q4 fund-holding detail rows.
# 2. Insert a new row 2 in "总计" (shifting the rest down) carrying the
#    2022-Q4 summary (count + market value), and append the row that is
#    pushed out the bottom (2020-Q4) as the new last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: build the new "2022-Q4" worksheet from a copy of "2022-Q3"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Column A (index) keeps the bordered/bold/centered style used throughout
# the workbook for that column - grab it from the template row above the
# rows we are about to add.
$q4.Range("A5").Copy()
$q4.Range("A6:A9").PasteSpecial(-4122)

# Columns B-G hold numeric-looking values that must stay TEXT (fund codes
# with leading zeros, percentages, etc.) - force text format before writing.
$q4.Range("B2:G9").NumberFormat = "@"

# -- row 2 --------------------------------------------------------------
$q4.Cells.Item(2,1).Value = 0
$q4.Cells.Item(2,2).Value = "009837"
$q4.Cells.Item(2,3).Value = "华夏磐锐一年定期开放混合A"
$q4.Cells.Item(2,4).Value = "14.15"
$q4.Cells.Item(2,5).Value = "75.21"
$q4.Cells.Item(2,6).Value = "3.82"
$q4.Cells.Item(2,7).Value = "0.5405"
$q4.Cells.Item(2,8).Value = 6

# -- row 3 --------------------------------------------------------------
$q4.Cells.Item(3,1).Value = 1
$q4.Cells.Item(3,2).Value = "005660"
$q4.Cells.Item(3,3).Value = "嘉实资源精选股票A"
$q4.Cells.Item(3,4).Value = "2.63"
$q4.Cells.Item(3,5).Value = "93.36"
$q4.Cells.Item(3,6).Value = "4.47"
$q4.Cells.Item(3,7).Value = "0.1176"
$q4.Cells.Item(3,8).Value = 8

# -- row 4 --------------------------------------------------------------
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "015443"
$q4.Cells.Item(4,3).Value = "惠升惠享启睿混合A"
$q4.Cells.Item(4,4).Value = "1.73"
$q4.Cells.Item(4,5).Value = "64.50"
$q4.Cells.Item(4,6).Value = "3.49"
$q4.Cells.Item(4,7).Value = "0.0604"
$q4.Cells.Item(4,8).Value = 8

# -- row 5 --------------------------------------------------------------
$q4.Cells.Item(5,1).Value = 3
$q4.Cells.Item(5,2).Value = "005661"
$q4.Cells.Item(5,3).Value = "嘉实资源精选股票C"
$q4.Cells.Item(5,4).Value = "1.06"
$q4.Cells.Item(5,5).Value = "93.36"
$q4.Cells.Item(5,6).Value = "4.47"
$q4.Cells.Item(5,7).Value = "0.0474"
$q4.Cells.Item(5,8).Value = 8

# -- row 6 --------------------------------------------------------------
$q4.Cells.Item(6,1).Value = 4
$q4.Cells.Item(6,2).Value = "005947"
$q4.Cells.Item(6,3).Value = "德邦民裕进取量化精选灵活配置混合A"
$q4.Cells.Item(6,4).Value = "0.50"
$q4.Cells.Item(6,5).Value = "90.73"
$q4.Cells.Item(6,6).Value = "7.34"
$q4.Cells.Item(6,7).Value = "0.0367"
$q4.Cells.Item(6,8).Value = 3

# -- row 7 --------------------------------------------------------------
$q4.Cells.Item(7,1).Value = 5
$q4.Cells.Item(7,2).Value = "009838"
$q4.Cells.Item(7,3).Value = "华夏磐锐一年定期开放混合C"
$q4.Cells.Item(7,4).Value = "0.39"
$q4.Cells.Item(7,5).Value = "75.21"
$q4.Cells.Item(7,6).Value = "3.82"
$q4.Cells.Item(7,7).Value = "0.0149"
$q4.Cells.Item(7,8).Value = 6

# -- row 8 --------------------------------------------------------------
$q4.Cells.Item(8,1).Value = 6
$q4.Cells.Item(8,2).Value = "005948"
$q4.Cells.Item(8,3).Value = "德邦民裕进取量化精选灵活配置混合C"
$q4.Cells.Item(8,4).Value = "0.09"
$q4.Cells.Item(8,5).Value = "90.73"
$q4.Cells.Item(8,6).Value = "7.34"
$q4.Cells.Item(8,7).Value = "0.0066"
$q4.Cells.Item(8,8).Value = 3

# -- row 9 --------------------------------------------------------------
$q4.Cells.Item(9,1).Value = 7
$q4.Cells.Item(9,2).Value = "015444"
$q4.Cells.Item(9,3).Value = "惠升惠享启睿混合C"
$q4.Cells.Item(9,4).Value = "0.00"
$q4.Cells.Item(9,5).Value = "64.50"
$q4.Cells.Item(9,6).Value = "3.49"
# G9 is a genuine number (0), unlike the other text-formatted G cells.
$q4.Range("G9").NumberFormat = "General"
$q4.Cells.Item(9,7).Value = 0
$q4.Cells.Item(9,8).Value = 8

# The "@" format was only needed to stop Excel from auto-converting the
# numeric-looking text while typing it in; switch the display format back
# to General now that the values are safely stored as text (this does not
# turn them back into numbers).
$q4.Range("B2:G9").NumberFormat = "General"

# ---------------------------------------------------------------------
# Step 2: shift "总计" down one row and insert the 2022-Q4 summary
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Restore column-A styling (border/bold/center) on the freshly inserted
# row, and drop the stray formatting Excel copied onto B2:D2 from the
# header row above.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 8
$total.Cells.Item(2,4).Value = 0.82

# The row-insert shifted the old rows down but left their 0-based index
# column (A) untouched - renumber it to stay sequential.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
$total.Cells.Item(7,1).Value = 5
$total.Cells.Item(8,1).Value = 6
$total.Cells.Item(9,1).Value = 7

# Keep the original active sheet/selection ("总计", A1) as it was before
# this edit.
$total.Activate()
[void]$total.Range("A1").Select()
